$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @(624, 0.796765087705063,   "Wilcoxon rank sum exact test",                       "two.sided", "p1 ~ media_v1"),
  @(523, 0.1633357758147733,  "Wilcoxon rank sum exact test",                       "two.sided", "p2 ~ media_v1"),
  @(812, 0.06440332494313779, "Wilcoxon rank sum exact test",                       "two.sided", "r1 ~ media_v1"),
  @(814, 0.06116922904381553, "Wilcoxon rank sum exact test",                       "two.sided", "r2 ~ media_v1"),
  @(820, 0.05225511575583854, "Wilcoxon rank sum exact test",                       "two.sided", "r3 ~ media_v1"),
  @(827, 0.04369004672153537, "Wilcoxon rank sum test with continuity correction",  "two.sided", "i1 ~ media_v1"),
  @(625, 0.8029677201588449,  "Wilcoxon rank sum test with continuity correction",  "two.sided", "i2 ~ media_v1"),
  @(537, 0.3256304515919216,  "Wilcoxon rank sum exact test",                       "two.sided", "p1 ~ media_v2"),
  @(555, 0.4373968125481567,  "Wilcoxon rank sum exact test",                       "two.sided", "p2 ~ media_v2"),
  @(733, 0.2122279251510636,  "Wilcoxon rank sum exact test",                       "two.sided", "r1 ~ media_v2"),
  @(740, 0.1841164168323889,  "Wilcoxon rank sum exact test",                       "two.sided", "r2 ~ media_v2"),
  @(742, 0.1766083506160379,  "Wilcoxon rank sum exact test",                       "two.sided", "r3 ~ media_v2"),
  @(743, 0.1717789922490576,  "Wilcoxon rank sum test with continuity correction",  "two.sided", "i1 ~ media_v2"),
  @(497, 0.1453250431908469,  "Wilcoxon rank sum test with continuity correction",  "two.sided", "i2 ~ media_v2"),
  @(393, 0.1117474076557653,  "Wilcoxon rank sum test with continuity correction",  "two.sided", "p1 ~ media_v3"),
  @(393, 0.1117474076557653,  "Wilcoxon rank sum test with continuity correction",  "two.sided", "p2 ~ media_v3"),
  @(668, 0.06368126808427531, "Wilcoxon rank sum test with continuity correction",  "two.sided", "r1 ~ media_v3"),
  @(678, 0.047689241617696,   "Wilcoxon rank sum test with continuity correction",  "two.sided", "r2 ~ media_v3"),
  @(699, 0.02482313013329353, "Wilcoxon rank sum test with continuity correction",  "two.sided", "r3 ~ media_v3"),
  @(615, 0.2347280906791234,  "Wilcoxon rank sum test with continuity correction",  "two.sided", "i1 ~ media_v3"),
  @(499, 0.7952970451071087,  "Wilcoxon rank sum test with continuity correction",  "two.sided", "i2 ~ media_v3")
)

$row = 2
foreach ($r in $data) {
  $ws.Cells.Item($row, 1).Value = $r[0]
  $ws.Cells.Item($row, 2).Value = $r[1]
  $ws.Cells.Item($row, 3).Value = $r[2]
  $ws.Cells.Item($row, 4).Value = $r[3]
  $ws.Cells.Item($row, 5).Value = $r[4]
  $row++
}
